$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header height shrinks from 60 to 30
$ws.Rows("1").RowHeight = 30

# Fill in the sales-percentage values for C2:C13 and normalize their
# formatting: thousands separator number format, wrapped + vertically
# centered text, default font/border (matches the other data columns).
$values = @{
    "C2"  = 6523
    "C3"  = 58504
    "C4"  = 132237
    "C5"  = 11612
    "C6"  = 57284
    "C7"  = 23644
    "C8"  = 66488
    "C9"  = 70022
    "C10" = 152621
    "C11" = 77549
    "C12" = 142997
    "C13" = 20379
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$col = $ws.Range("C2:C13")
$col.NumberFormat = "#,##0"
$col.WrapText = $true
$col.VerticalAlignment = -4108
$col.Borders.LineStyle = 1

# New (currently empty) formatted rows further down the sheet.
$lower = $ws.Range("B19:D31")
$lower.WrapText = $true
$lower.VerticalAlignment = -4108

$e31 = $ws.Range("E31")
$e31.NumberFormat = "#,##0"
$e31.WrapText = $true
$e31.VerticalAlignment = -4108

# Move the active selection to C15, matching where the author left off.
$ws.Range("C15").Select()
